$d = $word.ActiveDocument

# --- Change 1: remove the "Meta description: ..." paragraph that follows ---
# --- the H1 title at the top of the document.                           ---
[void]$d.Paragraphs(2).Range.Delete()

# --- Change 2: at the end of the document, replace the final (italic)   ---
# --- "Create a cartoon-style feature image..." paragraph with two new   ---
# --- paragraphs: a bold "Play Fairy Queen Free: Review and Gameplay"    ---
# --- paragraph, followed by an italic paragraph carrying the (moved)    ---
# --- meta-description sentence.                                         ---
$n = $d.Paragraphs.Count
[void]$d.Paragraphs($n).Range.Delete()

$dup = $d.Content.Duplicate
$dup.Collapse(0)

$xml = '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
  '<w:body>' + `
  '<w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Fairy Queen Free: Review and Gameplay</w:t></w:r></w:p>' + `
  '<w:p><w:r/><w:r><w:rPr><w:i/></w:rPr><w:t>Get ready to rule the fairy kingdom with Fairy Queen, an online slot game with free bonus spins and up to 9000x jackpot. Play now for free.</w:t></w:r></w:p>' + `
  '</w:body></w:document>'

[void]$dup.InsertXML($xml)
